# Update the "Fields" worksheet (sheet3):
#   - Add a "Field name=value" header in C1 and a "Field type" header in D1,
#     bold like the existing A1/B1 headers.
#   - Relabel the existing maintainer hyperlink cell (C2) so it reads
#     "maintainer=erik@ardoq.com" (keeping its mailto: hyperlink).
#   - Add "email" in D2 describing the type of the C2 field value.
#   - Leave the active selection on C2.
#
# The writes are ordered so that new shared-string table entries end up
# in the same order as the target workbook (email, Field type,
# maintainer=erik@ardoq.com, Field name=value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

$ws.Range("D2").Value = "email"
$ws.Range("D1").Value = "Field type"
$ws.Range("C2").Value = "maintainer=erik@ardoq.com"
$ws.Range("C1").Value = "Field name=value"

# Match the bold header formatting already used for A1/B1.
$ws.Range("C1:D1").Font.Bold = $true

# Reflect the selection state recorded in the workbook.
$ws.Range("C2").Select()
